$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.729.12"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "2.568.49"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'560.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.70%  "
$ws.Range("D6").Value = "'142.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("D9").Value = "2.575.25"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").Value = "'6.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("E12").Value = "  +7.93%  "
$ws.Range("D13").Value = "'0.340"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("D14").Value = "3.024.20"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").Value = "58.811.24"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").Value = "'21.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.81%  "
$ws.Range("E17").Value = "  +5.03%  "
$ws.Range("D18").Value = "2.580.14"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'4.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("D20").Value = "'334.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "'10.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("D22").Value = "'6.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'63.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  +6.27%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "'0.160"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("D28").Value = "'7.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.51%  "
$ws.Range("D29").Value = "0.0₃0777"
$ws.Range("E29").Value = "  +5.85%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("D32").Value = "'158.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.30%  "
$ws.Range("D33").Value = "'6.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("D34").Value = "'18.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("D35").Value = "'3.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.29%  "
$ws.Range("D36").Value = "'0.873"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.40%  "
$ws.Range("D37").Value = "'0.873"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.12%  "
$ws.Range("D38").Value = "'1.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.97%  "
$ws.Range("D39").Value = "'36.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("E40").Value = "  +4.38%  "
$ws.Range("D41").Value = "'290.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.82%  "
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  +3.01%  "
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").Value = "'10.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").Value = "'0.0531"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.11%  "
$ws.Range("D48").Value = "'19.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.92%  "
$ws.Range("D49").Value = "'123.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.13%  "
$ws.Range("D50").Value = "'0.0230"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").Value = "'18.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.22%  "
